# Generate Report for Handoff
# Adds a new tracked file (9ede6688-a181-487d-9709-b68a63b2ee1e.md) as a new
# row on each of the three worksheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ea03b7235da64b818b4b29b0ada181d1f2902b41"

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = "9ede6688-a181-487d-9709-b68a63b2ee1e.md"
$wsOverview.Range("B3").Value = "e2e\9ede6688-a181-487d-9709-b68a63b2ee1e.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-29 02:39:31"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "$repoBase/e2e/9ede6688-a181-487d-9709-b68a63b2ee1e.md", "", "", "e2e\9ede6688-a181-487d-9709-b68a63b2ee1e.md") | Out-Null
$wsOverview.Range("B3").Style = "HyperLink"

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A3").Value = "9ede6688-a181-487d-9709-b68a63b2ee1e.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = "9ede6688-a181-487d-9709-b68a63b2ee1e.bda54e574bff361b50215b95ca39336de05bfe11.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-29 02:39:27"
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("M3").Value = "'True"
$wsZhCn.Range("O3").Value = "'False"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "$repoBase/e2e/9ede6688-a181-487d-9709-b68a63b2ee1e.md", "", "", "9ede6688-a181-487d-9709-b68a63b2ee1e.md") | Out-Null
$wsZhCn.Range("A3").Style = "HyperLink"

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A3").Value = "9ede6688-a181-487d-9709-b68a63b2ee1e.md"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = "9ede6688-a181-487d-9709-b68a63b2ee1e.bda54e574bff361b50215b95ca39336de05bfe11.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-29 02:39:31"
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("M3").Value = "'True"
$wsDeDe.Range("O3").Value = "'False"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "$repoBase/e2e/9ede6688-a181-487d-9709-b68a63b2ee1e.md", "", "", "9ede6688-a181-487d-9709-b68a63b2ee1e.md") | Out-Null
$wsDeDe.Range("A3").Style = "HyperLink"
